# Updated symbol list on Mon Dec 26 15:48:50 UTC 2022 with GitHub Actions
# Refreshes the "Price" column (D) figures and two "Volume(1h)" (E) text
# labels that the scraper regenerated on this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells hold their numeric-looking price as TEXT (no leading
# apostrophe quirk, no quotePrefix) in the original workbook. Assigning a
# numeric-looking string straight to .Value lets Excel auto-coerce it to a
# real number, which would change the cell's stored type. Forcing the
# number format to Text ("@") before the assignment keeps it a string, and
# restoring the style afterwards (Style = "Normal") drops the now-unneeded
# explicit text format so the cell's style stays the same as before.
function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2"  "242.71"
Set-TextValue "D4"  "5.427"
Set-TextValue "D5"  "0.05888"
Set-TextValue "D8"  "0.8114"
Set-TextValue "D9"  "0.9586"
Set-TextValue "D11" "0.07444"
Set-TextValue "D12" "0.03268"
Set-TextValue "D13" "0.03066"
Set-TextValue "D14" "0.09347"
Set-TextValue "D15" "3.852"
Set-TextValue "D16" "0.001573"
Set-TextValue "D17" "0.04678"
Set-TextValue "D18" "0.0005929"
Set-TextValue "D19" "0.005888"
Set-TextValue "D20" "0.001254"
Set-TextValue "D21" "0.004898"
Set-TextValue "D22" "0.00006810"
Set-TextValue "D24" "2.128"
Set-TextValue "D26" "0.1309"
Set-TextValue "D27" "0.0002288"
Set-TextValue "D40" "0.03928"
Set-TextValue "D41" "0.006195"
Set-TextValue "D42" "0.1071"
Set-TextValue "D43" "0.003004"
Set-TextValue "D44" "0.009147"
Set-TextValue "D45" "0.00005216"
Set-TextValue "D47" "0.7312"
Set-TextValue "D48" "0.002298"

# Column E cells are plain (non-numeric-looking) text, so a direct .Value
# assignment keeps them as text without any extra formatting work.
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
